$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "Situação" (pass/fail) column
$ws.Range("Q1").Value = "Situação"

$green = 9498256   # RGB(144,238,144) -> FF90EE90
$red   = 8421616   # RGB(240,128,128) -> FFF08080

for ($r = 2; $r -le 160; $r++) {
    $pCell = $ws.Range("P$r")
    $pCell.Formula = "=AVERAGE(H$r,N$r,O$r)"

    $avg = $pCell.Value2

    $qCell = $ws.Range("Q$r")
    if ($avg -ge 10) {
        $qCell.Value = "Aprovado"
        $qCell.Interior.Color = $green
    } else {
        $qCell.Value = "Reprovado"
        $qCell.Interior.Color = $red
    }
}

$ws.Range("B2").Select()
